$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: write a value into a cell while preserving "text" storage for
# numeric-looking strings (the source workbook stores every data value, even
# numbers, as text in the shared-string table). Plain (non numeric-looking)
# text is safe to assign directly.
function Set-TextValue {
    param($addr, $val)

    $isNumeric = $val -match '^-?[0-9]+(\.[0-9]+)?$'

    $rng = $ws.Range($addr)
    if ($isNumeric) {
        $rng.NumberFormat = "@"
        $rng.Value = $val
        $rng.Style = "Normal"
    } else {
        $rng.Value = $val
    }
}

# ------------------------------------------------------------------
# 1) Update the "Pertuis d'Antioche" row (row 8) measurements
# ------------------------------------------------------------------
Set-TextValue "E8" "18.27"
Set-TextValue "F8" "57.87"
Set-TextValue "G8" "103.86"
Set-TextValue "I8" "301.97"
Set-TextValue "J8" "111.63"

# ------------------------------------------------------------------
# 2) Swap the contents of row 22 (Sainte-Mère-Église) and row 23
#    (Massif des Maures) across columns D:G and I:K (H is identical
#    in both rows, so it is left untouched).
# ------------------------------------------------------------------
Set-TextValue "D22" "Massif des Maures"
Set-TextValue "E22" "11.79"
Set-TextValue "F22" "23.05"
Set-TextValue "G22" "145.16"
Set-TextValue "I22" "540.52"
Set-TextValue "J22" "282.41"
Set-TextValue "K22" "CW"

Set-TextValue "D23" "Sainte-Mère-Église"
Set-TextValue "E23" "177.23"
Set-TextValue "F23" "0.9"
Set-TextValue "G23" "1.87"
Set-TextValue "I23" "378.14"
Set-TextValue "J23" "1166.87"
Set-TextValue "K23" "CCW"

# ------------------------------------------------------------------
# 3) Replace the two custom column widths (C and D) with a single
#    uniform width applied to columns A through K.
# ------------------------------------------------------------------
$ws.Range("A:K").ColumnWidth = 16.1
